# update rg+curseurs pour ssrha gme pypmsi
#
# Appends 7 new rows (182-188) to Feuil1 describing the new "ssrha_gme"
# table / "zgme" column regex-cursor entries (years 2017-2023), matching
# the shape/style of every other table block already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 7 fresh rows right after the last existing row (181), cloned from
# row 181 itself so the new rows pick up the same cell styles (s="2" on
# column A, s="1" on column D) used throughout the rest of the table.
for ($i = 0; $i -lt 7; $i++) {
  $ws.Rows.Item(181).Copy()
  $ws.Rows.Item(182).Insert()
}

# table, an (year), z, rg (regex), curseur (regex length)
$newRows = @(
  @(182, 2017, ".{1,13}", 13),
  @(183, 2018, ".{1,13}", 13),
  @(184, 2019, ".{1,13}", 13),
  @(185, 2020, ".{1,13}", 13),
  @(186, 2021, ".{1,13}", 13),
  @(187, 2022, ".{1,14}", 14),
  @(188, 2023, ".{1,14}", 14)
)

foreach ($item in $newRows) {
  $r     = $item[0]
  $year  = $item[1]
  $regex = $item[2]
  $curs  = $item[3]

  $ws.Cells.Item($r, 1).Value = "ssrha_gme"
  $ws.Cells.Item($r, 2).Value = $year
  $ws.Cells.Item($r, 3).Value = "zgme"
  $ws.Cells.Item($r, 4).Value = $regex
  $ws.Cells.Item($r, 5).Value = $curs
}

# Match the author's final cursor position/selection
$ws.Range("B189").Select()
